$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "303.64"
Set-TextValue $ws.Range("E2") "5.92%"

Set-TextValue $ws.Range("D3") "31.80"
Set-TextValue $ws.Range("E3") "8.04%"

Set-TextValue $ws.Range("D4") "5.223"
Set-TextValue $ws.Range("E4") "2.20%"

Set-TextValue $ws.Range("D5") "0.07265"
Set-TextValue $ws.Range("E5") "8.27%"

Set-TextValue $ws.Range("D6") "7.792"
Set-TextValue $ws.Range("E6") "6.09%"

Set-TextValue $ws.Range("D7") "3.750"
Set-TextValue $ws.Range("E7") "8.90%"

Set-TextValue $ws.Range("D8") "1.444"
Set-TextValue $ws.Range("E8") "3.94%"

Set-TextValue $ws.Range("D9") "0.9064"
Set-TextValue $ws.Range("E9") "-0.97%"

Set-TextValue $ws.Range("D10") "0.01630"
Set-TextValue $ws.Range("E10") "2,412.89%"

Set-TextValue $ws.Range("D11") "0.1675"
Set-TextValue $ws.Range("E11") "5.54%"

Set-TextValue $ws.Range("D12") "0.07475"
Set-TextValue $ws.Range("E12") "9.81%"

Set-TextValue $ws.Range("D13") "0.07915"
Set-TextValue $ws.Range("E13") "2.87%"

Set-TextValue $ws.Range("D14") "0.02980"
Set-TextValue $ws.Range("E14") "1.48%"

Set-TextValue $ws.Range("D15") "0.09939"
Set-TextValue $ws.Range("E15") "10.58%"

Set-TextValue $ws.Range("D16") "0.001505"
Set-TextValue $ws.Range("E16") "-5.44%"

Set-TextValue $ws.Range("D17") "0.04536"
Set-TextValue $ws.Range("E17") "1.33%"

Set-TextValue $ws.Range("D18") "0.006309"
Set-TextValue $ws.Range("E18") "0.51%"

Set-TextValue $ws.Range("D19") "3.471"
Set-TextValue $ws.Range("E19") "0.58%"

Set-TextValue $ws.Range("D20") "2.224"
Set-TextValue $ws.Range("E20") "-0.22%"

Set-TextValue $ws.Range("E21") "4.26%"

Set-TextValue $ws.Range("D22") "0.1334"
Set-TextValue $ws.Range("E22") "1.70%"

Set-TextValue $ws.Range("D23") "4.279"
Set-TextValue $ws.Range("E23") "5.06%"

Set-TextValue $ws.Range("D25") "0.001225"
Set-TextValue $ws.Range("E25") "2.44%"

Set-TextValue $ws.Range("D26") "0.004410"
Set-TextValue $ws.Range("E26") "7.11%"

Set-TextValue $ws.Range("D27") "0.0001305"
Set-TextValue $ws.Range("E27") "8.88%"

Set-TextValue $ws.Range("D28") "0.0001753"
Set-TextValue $ws.Range("E28") "8.30%"

Set-TextValue $ws.Range("D40") "0.04478"
Set-TextValue $ws.Range("E40") "4.99%"

Set-TextValue $ws.Range("D41") "0.007204"
Set-TextValue $ws.Range("E41") "7.07%"

Set-TextValue $ws.Range("D42") "0.1340"
Set-TextValue $ws.Range("E42") "8.07%"

Set-TextValue $ws.Range("D43") "0.002349"
Set-TextValue $ws.Range("E43") "4.99%"

Set-TextValue $ws.Range("D44") "0.01280"
Set-TextValue $ws.Range("E44") "7.22%"

Set-TextValue $ws.Range("D45") "0.00006089"
Set-TextValue $ws.Range("E45") "6.70%"

Set-TextValue $ws.Range("E46") "-3.44%"

Set-TextValue $ws.Range("D47") "0.01618"
Set-TextValue $ws.Range("E47") "7.34%"
